$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had May 2020 weekday entries on every other row (2,4,6,...42)
# with blank rows in between. The new sheet has June 2020 weekday entries
# written consecutively on rows 2-23 with no gaps. Clear the old range first.
$ws.Range("A2:B42").ClearContents()

# Make sure the date-looking text stays as literal text (not auto-converted
# to a date serial number) when assigned.
$ws.Range("A2:A23").NumberFormat = "@"

# New data: June 2020 weekdays, written consecutively starting at row 2.
$dates = @(
  "06/01/20","06/02/20","06/03/20","06/04/20","06/05/20",
  "06/08/20","06/09/20","06/10/20","06/11/20","06/12/20",
  "06/15/20","06/16/20","06/17/20","06/18/20","06/19/20",
  "06/22/20","06/23/20","06/24/20","06/25/20","06/26/20",
  "06/29/20","06/30/20"
)

$row = 2
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 2).Value = "June"
    $row = $row + 1
}

# Update the active-cell selection to match the new state.
$ws.Range("H16").Select()
